# Weekly update: a new price record was added to the dataset.
# It was inserted as a new row at sheet row 202, pushing the existing
# rows 202..231 down to 203..232 (dimension grows from A1:R231 to A1:R232).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 202 (this shifts rows 202-231 down
# to 203-232 and automatically extends the used range / dimension).
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(202, 1).Value = 4
$ws.Cells.Item(202, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(202, 3).Value = "Los Lagos"
$ws.Cells.Item(202, 4).Value = 44491
$ws.Cells.Item(202, 5).Value = 10
$ws.Cells.Item(202, 6).Value = 100114013
$ws.Cells.Item(202, 7).Value = "Zanahoria"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 750
$ws.Cells.Item(202, 11).Value = 12000
$ws.Cells.Item(202, 12).Value = 12000
$ws.Cells.Item(202, 13).Value = 12000
$ws.Cells.Item(202, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(202, 15).Value = "Región de Ñuble"
$ws.Cells.Item(202, 16).Value = 600
$ws.Cells.Item(202, 17).Value = 20
$ws.Cells.Item(202, 18).Value = "Hortaliza"
